$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the CODE row entries for "linkedlist" work
$ws.Range("D6").Value = "1.linkedlist"
$ws.Range("C7").Value = "2. linkedlist"
$ws.Range("F6").Value = "1.Created linkedlist"

# Update the active selection to C16, matching the saved view state
$ws.Range("C16").Select()
